$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change B2 and C2 from numeric hour values to text time strings
$ws.Range("B2").Value = "19:30"
$ws.Range("C2").Value = "21:00"
